# Replace the 11-row (A:B) MOC data with the updated 31-row data set and
# add a third column (C) of CSV-exported values (all zero).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = @(
    0, 22.473693933007748, 29.722275188127529, 34.153532513107017, 38.261896087222766,
    42.268384043729085, 46.279251758217093, 50.357539951497323, 54.546543300144471, 58.879476474951218,
    63.384097679175881, 68.085177075997294, 73.005936021122082, 78.168952826997582, 83.596774571649362,
    89.312359600325365, 95.339419730144414, 101.70270247624951, 108.4282380431073, 115.54356700780011,
    123.07795946129629, 131.06263327736957, 139.53097729473504, 148.51878405240365, 158.06449603956796,
    168.20946904965939, 178.99825606763144, 190.47891511183616, 202.70334456184622, 215.72764971047636,
    229.61254457070822
)

$colB = @(
    15.227067222484973, 25.199196518943083, 28.295530262795801, 30.079942447862663, 31.668230001666924,
    33.153362471893381, 34.576934123207032, 35.960820607690735, 37.317548943874876, 38.654543003430874,
    39.976125065687171, 41.284553754015711, 42.580597778855143, 43.863865204006082, 45.132994077560035,
    46.385759120780676, 47.619124318470661, 48.829258371390218, 50.011522909929695, 51.160439299863903,
    52.269637413571225, 53.331788182791115, 54.338520709398907, 55.280323974474229, 56.146432628277381,
    56.924695889012611, 57.601428178673331, 58.161239748114831, 58.586845168418044, 58.858847174625488,
    58.95549292705909
)

$rowCount = $colA.Count

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    $ws.Cells.Item($row, 3).Value = 0
}

# Scroll/selection state left by the author after entering the new data.
$ws.Range("E24").Select()
